# "times update with calculator hours"
# Adds a small "Days Calculator" block (rows 33-35) below the existing
# hours summary, widens/normalises the data-entry columns, marks a few
# more "Mondag" cells with the legend colours, and applies a one-decimal
# number format to the summary numbers in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: make A:F a uniform width -------------------------
$ws.Columns("A:F").ColumnWidth = 17.140625

# --- Re-colour the "Mondag" (column B) cells for weeks 5-25 ----------
# Reuse the existing "School" colouring (as already used on B2/B4) and
# the existing "At intership" colouring (as already used on B5) by
# copying their formats onto the remaining week rows, rather than
# inventing new styles.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B6:B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B11:B26").PasteSpecial(-4122) | Out-Null

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- One-decimal number format for the hours/days summary ------------
$ws.Range("F28:F33").NumberFormat = "0.0"

# --- New "Days Calculator" block (rows 33-35) -------------------------
$ws.Range("A33").Value = "Days Calculator (H)"
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = "One day (H)"
$ws.Range("B34").Formula = "=C30"
$ws.Range("A35").Value = "Answer: (Day)"
$ws.Range("B35").Formula = "=B33/B34"

# --- Move the active selection the way the author left it ------------
$ws.Range("A36").Select() | Out-Null
